# Remove the "Delete button is now working for first to do list item but
# not for later items." list-paragraph entirely (including its paragraph
# mark), leaving the surrounding paragraphs untouched.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Delete button is now working for first to do list item but not for later items.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
